# Update cryptos list (GitHub Actions refresh) — apply latest price/volume
# snapshot to the existing worksheet, plus a ranking swap between
# Avalanche and Chainlink (rows 22 and 23 trade places in the list while
# keeping their rank numbers in column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22 / 23: Avalanche and Chainlink swap positions -------------------
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.18"
$ws.Range("E22").Value = "  +2.98%  "

$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.40"
$ws.Range("E23").Value = "  +0.99%  "

# --- Price (D) / Volume(1h) (E) refresh for all other rows -----------------
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.675.57"
$ws.Range("E2").Value = "  +1.66%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.635.64"
$ws.Range("E3").Value = "  +1.94%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.92"
$ws.Range("E5").Value = "  +0.14%  "

$ws.Range("E6").Value = "  +1.64%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  +1.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0623"
$ws.Range("E9").Value = "  +1.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.04"
$ws.Range("E10").Value = "  +3.70%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0840"
$ws.Range("E11").Value = "  +3.52%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.863.62"
$ws.Range("E12").Value = "  +1.95%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.639.62"
$ws.Range("E13").Value = "  +2.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.07"
$ws.Range("E14").Value = "  +1.42%  "

$ws.Range("E15").Value = "  +2.76%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.670.30"
$ws.Range("E16").Value = "  +1.77%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.00"
$ws.Range("E17").Value = "  +1.65%  "

$ws.Range("E18").Value = "  +1.64%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "209.05"
$ws.Range("E19").Value = "  +4.39%  "

$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("E21").Value = "  +1.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.72"
$ws.Range("E25").Value = "  +1.79%  "

$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("E27").Value = "  -0.74%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.76"
$ws.Range("E28").Value = "  +3.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.37"
$ws.Range("E29").Value = "  +1.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0519"
$ws.Range("E30").Value = "  +5.62%  "

$ws.Range("E31").Value = "  -0.35%  "

$ws.Range("E32").Value = "  +1.19%  "

$ws.Range("E33").Value = "  +0.64%  "

$ws.Range("E34").Value = "  +0.79%  "

$ws.Range("E35").Value = "  -0.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.169.33"

$ws.Range("E37").Value = "  -1.42%  "

$ws.Range("E38").Value = "  +3.14%  "

$ws.Range("E39").Value = "  +0.02%  "

$ws.Range("E40").Value = "  +1.37%  "

$ws.Range("E41").Value = "  +0.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.793"
$ws.Range("E42").Value = "  +1.18%  "

$ws.Range("E43").Value = "  +0.77%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.775.55"
$ws.Range("E44").Value = "  +2.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.05"
$ws.Range("E45").Value = "  +0.12%  "

$ws.Range("E46").Value = "  +1.01%  "

$ws.Range("E47").Value = "  -2.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.67"
$ws.Range("E48").Value = "  +1.25%  "

$ws.Range("E49").Value = "  +1.64%  "

$ws.Range("E50").Value = "  +0.55%  "

$ws.Range("E51").Value = "  +4.44%  "
